# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that hold the "Table_0" styled grids) get
#    switched from the deck's custom table style to PowerPoint's built-in
#    "No Style, No Grid" table style.
# 2) The presentation's theme colour scheme is switched from the "Integral /
#    Red Violet" palette back to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Retarget the table styles -------------------------------------------------
$builtInStyleId = "{4CB49EAC-77D8-4628-A45C-07985E5DA99B}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($builtInStyleId)
        }
    }
}

# --- 2) Restore the stock "Office Theme" colour scheme ----------------------------
function Get-ComRGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (Office Theme defaults)
$officeThemeColors = @(
    (Get-ComRGB 0x00 0x00 0x00),  # dk1      000000
    (Get-ComRGB 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (Get-ComRGB 0x44 0x54 0x6A),  # dk2      44546A
    (Get-ComRGB 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (Get-ComRGB 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (Get-ComRGB 0xED 0x7D 0x31),  # accent2  ED7D31
    (Get-ComRGB 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (Get-ComRGB 0xFF 0xC0 0x00),  # accent4  FFC000
    (Get-ComRGB 0x44 0x72 0xC4),  # accent5  4472C4
    (Get-ComRGB 0x70 0xAD 0x47),  # accent6  70AD47
    (Get-ComRGB 0x05 0x63 0xC1),  # hlink    0563C1
    (Get-ComRGB 0x95 0x4F 0x72)   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
